$wb = $excel.ActiveWorkbook

# --- Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B6").Value = -13351.53190127437
$ws.Range("B7").Value = 10922087.44374209
$ws.Range("B8").Value = 25197720.06796782
$ws.Range("B10").Value = 2688784.803334876

# --- Costs and Revenues ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("E2").Value = 61484.09347675643
$ws.Range("H2").Value = 64565.43974857162
$ws.Range("I2").Value = 55848.12409846896
$ws.Range("J2").Value = 65698.51352618316
$ws.Range("K2").Value = 43379.07225362272
$ws.Range("L2").Value = 48995.05304514673
$ws.Range("M2").Value = 62448.52877130933
$ws.Range("N2").Value = 45450.49222258999
$ws.Range("O2").Value = 50696.26030812534
$ws.Range("P2").Value = 37157.69159838425
$ws.Range("E3").Value = 133100.0000000001
$ws.Range("E4").Value = 28586.87318168494
$ws.Range("G4").Value = 5114.001716423241
$ws.Range("H4").Value = 31668.21945350014
$ws.Range("I4").Value = 22950.90380339748
$ws.Range("J4").Value = 32801.29323111167
$ws.Range("K4").Value = 10481.85195855123
$ws.Range("L4").Value = 16097.83275007525
$ws.Range("M4").Value = 29551.30847623785
$ws.Range("N4").Value = 12553.2719275185
$ws.Range("O4").Value = 17799.04001305386
$ws.Range("P4").Value = 4260.471303312761
$ws.Range("E6").Value = -103814.3918099369
$ws.Range("H6").Value = 29285.60819006316
$ws.Range("I6").Value = 29285.60819006316
$ws.Range("K6").Value = 29285.60819006316
$ws.Range("L6").Value = 29285.60819006316
$ws.Range("M6").Value = 29285.60819006315
$ws.Range("N6").Value = 29285.60819006316
$ws.Range("O6").Value = 29285.60819006316
$ws.Range("P6").Value = 29285.60819006316

# --- Fed-in Capacity ---
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("J20").Value = 124.5190384721106
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 113.4004983079896
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 117.8828208804077
$ws.Range("Q20").Value = 150.3839754851235
$ws.Range("M21").Value = 51.84373129681028
$ws.Range("O22").Value = 96.22962838366004
$ws.Range("P22").Value = 101.5955875616828
$ws.Range("K23").Value = 0
$ws.Range("Q23").Value = 150.3839754851235
$ws.Range("K24").Value = 80.29914934735042
$ws.Range("O24").Value = 57.81213424001893
$ws.Range("P24").Value = 65.92768427608706
$ws.Range("J26").Value = 0
$ws.Range("O26").Value = 117.8828208804077
$ws.Range("P26").Value = 135.4597561231036
$ws.Range("Q26").Value = 0
$ws.Range("J27").Value = 93.17061249236157
$ws.Range("O27").Value = 57.81213424001893
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 92.09541281912071
$ws.Range("N28").Value = 81.96869489115805
$ws.Range("O28").Value = 96.22962838366004
$ws.Range("P28").Value = 0
$ws.Range("R29").Value = 0
$ws.Range("K30").Value = 80.29914934735042
$ws.Range("P30").Value = 65.92768427608706
$ws.Range("Q30").Value = 94.49434172313325
$ws.Range("L32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("O34").Value = 96.22962838366004
$ws.Range("P34").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = 110.5750244233121
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("J38").Value = 124.5190384721106
$ws.Range("Q39").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("Q41").Value = 150.3839754851235
$ws.Range("R41").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("J46").Value = 0

# --- Unmet Demand ---
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("L11").Value = 130.6648563030561
$ws.Range("M11").Value = 113.4004983079896
$ws.Range("O11").Value = 117.8828208804077
$ws.Range("K12").Value = 80.29914934735042
$ws.Range("L12").Value = 61.18167021676314
$ws.Range("M12").Value = 51.84373129681028
$ws.Range("O12").Value = 57.81213424001893
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 135.370731907559
$ws.Range("L20").Value = 130.6648563030561
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 110.5750244233121
$ws.Range("O20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("K23").Value = 135.370731907559
$ws.Range("Q23").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("J26").Value = 124.5190384721106
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 150.3839754851235
$ws.Range("J27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("L28").Value = 90.4687457914608
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 101.5955875616828
$ws.Range("R29").Value = 173.7492132756177
$ws.Range("K30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("L32").Value = 130.6648563030561
$ws.Range("O32").Value = 117.8828208804077
$ws.Range("P32").Value = 135.4597561231036
$ws.Range("Q32").Value = 150.3839754851235
$ws.Range("K33").Value = 80.29914934735042
$ws.Range("L33").Value = 61.18167021676314
$ws.Range("M33").Value = 51.84373129681028
$ws.Range("N33").Value = 38.66169381481656
$ws.Range("O33").Value = 57.81213424001893
$ws.Range("P33").Value = 65.92768427608706
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 101.5955875616828
$ws.Range("K35").Value = 135.370731907559
$ws.Range("L35").Value = 130.6648563030561
$ws.Range("N35").Value = 0
$ws.Range("M36").Value = 51.84373129681028
$ws.Range("N36").Value = 38.66169381481656
$ws.Range("O36").Value = 57.81213424001893
$ws.Range("J38").Value = 0
$ws.Range("Q39").Value = 94.49434172313325
$ws.Range("N40").Value = 81.96869489115805
$ws.Range("K41").Value = 135.370731907559
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = 173.7492132756177
$ws.Range("J42").Value = 93.17061249236157
$ws.Range("K42").Value = 80.29914934735042
$ws.Range("L42").Value = 61.18167021676314
$ws.Range("P42").Value = 65.92768427608706
$ws.Range("J44").Value = 124.5190384721106
$ws.Range("J46").Value = 105.873818686614

# --- Household Surplus ---
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B5").Value = 298170.9753490015
$ws.Range("B8").Value = 314762.8398895449
$ws.Range("B9").Value = 267823.4479274536
$ws.Range("B10").Value = 320864.0063843761
$ws.Range("B11").Value = 200682.3995321278
$ws.Range("B12").Value = 230922.2961018724
$ws.Range("B13").Value = 303364.0884735172
$ws.Range("B14").Value = 211836.1993650283
$ws.Range("B15").Value = 240082.6429025264
$ws.Range("B16").Value = 167182.6575423821
